$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (price + 1h volume change) per upstream diff.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 45/46 swap Monero/Bittensor order.

# Row 2
$ws.Range("D2").Value = "61.501.84"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3
$ws.Range("D3").Value = "2.930.81"
$ws.Range("E3").Value = "  +0.72%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.41"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.67"
$ws.Range("E6").Value = "  -1.07%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  -1.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.95"
$ws.Range("E9").Value = "  +0.91%  "

# Row 11
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000223"
$ws.Range("E12").Value = "  -0.63%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.19"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14
$ws.Range("E14").Value = "  +0.73%  "

# Row 15
$ws.Range("D15").Value = "3.415.95"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16
$ws.Range("D16").Value = "61.457.11"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
$ws.Range("D17").Value = "2.930.01"
$ws.Range("E17").Value = "  +0.63%  "

# Row 18
$ws.Range("E18").Value = "  -0.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.34"
$ws.Range("E19").Value = "  +0.97%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.55"
$ws.Range("E20").Value = "  +1.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.671"
$ws.Range("E21").Value = "  -0.56%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.40"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -1.04%  "

# Row 25
$ws.Range("E25").Value = "  -0.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -1.14%  "

# Row 28
$ws.Range("E28").Value = "  -3.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("E31").Value = "  +1.08%  "

# Row 32
$ws.Range("E32").Value = "  +1.99%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0876"
$ws.Range("E34").Value = "  +3.20%  "

# Row 35
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("E37").Value = "  -0.76%  "

# Row 38
$ws.Range("E38").Value = "  +0.54%  "

# Row 39
$ws.Range("E39").Value = "  +0.43%  "

# Row 40
$ws.Range("E40").Value = "  -0.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.94"
$ws.Range("E41").Value = "  +5.26%  "

# Row 42
$ws.Range("E42").Value = "  -2.09%  "

# Row 43
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("D44").Value = "2.694.41"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "363.45"
$ws.Range("E45").Value = "  -2.88%  "

# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.00"
$ws.Range("E46").Value = "  +0.73%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.51"
$ws.Range("E48").Value = "  -0.66%  "

# Row 49
$ws.Range("E49").Value = "  -1.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -0.27%  "

# Row 51
$ws.Range("E51").Value = "  +0.74%  "
